$d = $word.ActiveDocument

# 1) Merge the "Stubhub" split run back into the surrounding text so the
#    proofErr-wrapped "Stubhub" run disappears and the sentence becomes a
#    single contiguous run of text.
$r = $d.Content.Find.Execute(
    "A link to the artists page and a calendar on the artists" + [char]8217 + " page that shows upcoming performances and venues with links to ticket vendors like Stubhub and Ticketmaster  ",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

# 2) Append the missing paragraph text after "Maple Hill "
$rng = $d.Content
$found = $rng.Find.Execute("Maple Hill ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter("has provided some useful information about generally likes in ANY website. This might be helpful in getting in vague understanding of what a musician likes in terms of usability and design. For example, it" + [char]8217 + "s easy to imagine that an artist might like a website that has lot of color and artsy design. Yet the guitarist, who asked to be called Cyclops, preferred a simple site without a lot of buttons and distractions that really gets in the way of the information that he needs. With this in mind, the need to appeal to the performers visual senses is less important than providing a clear, straightforward path to providing them with the useful information about the festival.")
}
